$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.413.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '2.279.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''303.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''95.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.29%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.502'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.24%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.493'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.64%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''34.99'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.86%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''18.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''6.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.09%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '2.634.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.14%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '2.271.73'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''0.769'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '42.331.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''12.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '0.0₃0886'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.65%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''5.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.94%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''67.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.00%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''235.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.79%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '  +1.12%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''2.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.16%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''24.60'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.26%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '  +17.00%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''167.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''8.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.90%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''32.21'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.48%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''4.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.04%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''17.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''4.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.14%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '  -2.38%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '  -2.49%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''2.64'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = '1.981.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.0274'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.05%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''9.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.44%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '  +2.12%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''2.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.47%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '  -2.72%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''2.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.52%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''53.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.89%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '2.502.50'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''70.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.51%  '
$ws.Range("E51").Style = "Normal"
